# Add a 2020 data column (column N) to the SDG indicator table, mirroring
# the formatting of the existing 2019 column (column M), then select M25
# as the last active cell (matching the saved view state of the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that receive a new value in column N, keyed by row number.
$newValues = [ordered]@{
    3  = 2020
    5  = 2198.6999999999998
    6  = 132.69999999999999
    7  = 242.9
    8  = 203.3
    9  = 202.8
    10 = 284.7
    11 = 294.89999999999998
    12 = 802.5
    13 = 28.1
    14 = 6.8
    16 = 27.4
    17 = 17.5
    18 = 24.7
    19 = 31.5
    20 = 30.4
    21 = 24.8
    22 = 30.7
    23 = 30.1
    24 = 21.2
    25 = 11.6
}

foreach ($row in $newValues.Keys) {
    $src = $ws.Cells.Item($row, 13)   # column M
    $dst = $ws.Cells.Item($row, 14)   # column N
    # Copy formatting (and the style index) from the 2019 (M) cell onto the
    # new 2020 (N) cell, then overwrite with the correct 2020 value.
    $src.Copy($dst)
    $dst.Value = $newValues[$row]
}

# Row 15 is a blank separator row; column N just needs to pick up the same
# (empty) style as column M, with no value.
$ws.Cells.Item(15, 13).Copy($ws.Cells.Item(15, 14))

# Restore the selection that was active when the workbook was last saved.
$ws.Range("M25").Select()
